# aggiornamento fino a 27/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$rows = @(
    @(44330, 0, 2, 218.3406113537118),
    @(44331, 0, 2, 218.3406113537118),
    @(44332, 0, 2, 218.3406113537118),
    @(44333, 0, 1, 109.1703056768559),
    @(44334, 0, 1, 109.1703056768559),
    @(44335, 0, 1, 109.1703056768559),
    @(44336, 0, 0, 0),
    @(44337, 0, 0, 0),
    @(44338, 0, 0, 0),
    @(44339, 0, 0, 0),
    @(44340, 0, 0, 0),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$startRow = 256

# The date cell in column A (e.g. A255) carries a custom style (border,
# bold, centered, date number format). Copy that formatting down onto
# each new date cell instead of Range.Style (which hands back a raw
# COM wrapper that doesn't assign cleanly here).
$templateCell = $ws.Cells.Item($startRow - 1, 1)
$templateCell.Copy() | Out-Null

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $cellA.Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
